$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1) Delete the "Unit" sheet
$wsUnit = $wb.Worksheets.Item("Unit")
$wsUnit.Delete()

# 2) Update the "@prefix" sheet: row16 'unit' -> 'unitLength', with new URL
$wsPrefix = $wb.Worksheets.Item("@prefix")
$wsPrefix.Range("A16").Value = "unitLength"
$wsPrefix.Range("B16").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/Unit/UnitLength#"

# 3) Update the "size" sheet: D5 'unit:um' -> 'unitLength:um'
$wsSize = $wb.Worksheets.Item("size")
$wsSize.Range("D5").Value = "unitLength:um"

Write-Host "done"
